$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10: becomes the "/api/rooms/:id" row. Also fix D10's leftover
#     "s=0" style (a pre-existing quirk) to match the rest of the row (s=2),
#     by copying the format from a same-row neighbour that already has it ---
$ws.Range("A10").Value = "/api/rooms/:id"
$ws.Range("C10").Value = "{data:{room:{}}}"
$ws.Range("E10").Copy()
$ws.Range("D10").PasteSpecial(-4122)

# --- Row 11: new "/api/rooms/create" row (previously row 10's "create" row has been
#     pushed down because a new "/api/rooms/:id" row was inserted above it at row 10).
#     D11 also had the stray "s=0" style, so fix it the same way before writing its value ---
$ws.Range("A11").Value = "/api/rooms/create"
$ws.Range("B11").Value = "{number,baseRent}"
$ws.Range("C11").Value = "{data:{room:{}}}"
$ws.Range("E11").Copy()
$ws.Range("D11").PasteSpecial(-4122)
$ws.Range("D11").Value = "The newly created Room"

# --- Row 12: cleared out (the old "/api/tenants" + "{data:{rooms:[]}}" row content
#     moved down) ---
$ws.Range("A12").Value = ""
$ws.Range("B12").Value = ""
$ws.Range("C12").Value = ""
$ws.Range("D12").Value = ""

# --- Row 13: "/api/tenants" row, now with its own proper return type ---
$ws.Range("A13").Value = "/api/tenants"
$ws.Range("C13").Value = "{data:{tenants:[]}}"

# --- Row 14: new "/api/tenants/:id" row ---
$ws.Range("A14").Value = "/api/tenants/:id"
$ws.Range("C14").Value = "{data:{tenant:{}}}"

# --- Row 15: new "/api/tenants/create" row with full detail ---
$ws.Range("A15").Value = "/api/tenants/create"
$ws.Range("B15").Value = "{name,phoneNumber,aadharCard,room}"
$ws.Range("C15").Value = "{data:{tenant:{}}}"
$ws.Range("D15").Value = "The newly created Tenant"

# --- Row 17: new "/api/transactions/?room" row. D17 had the same stray
#     "s=0" style quirk even though it stays blank, so fix it too ---
$ws.Range("A17").Value = "/api/transactions/?room"
$ws.Range("C17").Value = "{data:{transactions:[]}}"
$ws.Range("E17").Copy()
$ws.Range("D17").PasteSpecial(-4122)

# --- Row 18: new "/api/transactions/:id" row ---
$ws.Range("A18").Value = "/api/transactions/:id"
$ws.Range("C18").Value = "{data:{transaction:{}}}"

# --- Row 19: new "/api/tenants/create/?roomNumber" row with full detail ---
$ws.Range("A19").Value = "/api/tenants/create/?roomNumber"
$ws.Range("B19").Value = "{room,balance,transfer,remarks}"
$ws.Range("C19").Value = "{data:{transaction:{}}}"
$ws.Range("D19").Value = "The newly created Transaction"

# --- Widen columns A and B (target display widths 37.64 / 41.68 chars;
#     Excel COM snaps ColumnWidth to whole-pixel increments, so the closest
#     achievable COM input is used) ---
$ws.Columns.Item(1).ColumnWidth = 36.8
$ws.Columns.Item(2).ColumnWidth = 40.85

# --- Extend the used range with two more (blank) rows, 33 and 34, formatted
#     the same as the rest of the table body (row 32) ---
for ($col = 1; $col -le 35; $col++) {
    $ws.Cells.Item(33, $col).Value = "x"
    $ws.Cells.Item(34, $col).Value = "x"
}
$ws.Range("A32:AI32").Copy()
$ws.Range("A33:AI34").PasteSpecial(-4122)
for ($col = 1; $col -le 35; $col++) {
    $ws.Cells.Item(33, $col).Value = ""
    $ws.Cells.Item(34, $col).Value = ""
}

# --- Update active selection ---
$ws.Range("B12").Select()

$wb.Save()
